# "Generate Report for Archive"
#
# The localization status moved from "Ready for handoff" to "In Translation".
# That status string is shown in four places:
#   - Overview!E2  (zh-cn status)
#   - Overview!F2  (de-de status)
#   - zh-cn!C2     (Status column)
#   - de-de!C2     (Status column)
# Updating the text also shrinks the "Status" columns, so their widths are
# re-fit to match the shorter value.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
$newColumnWidth = 12.576851254417766

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Sanity-check the cells hold the expected "before" text, then replace it.
$targets = @(
    @{ Sheet = $overview; Cell = "E2" },
    @{ Sheet = $overview; Cell = "F2" },
    @{ Sheet = $zhcn;     Cell = "C2" },
    @{ Sheet = $dede;     Cell = "C2" }
)

foreach ($target in $targets) {
    $range = $target.Sheet.Range($target.Cell)
    $current = $range.Value()
    if (($current -is [string]) -and ($current -eq $oldStatus)) {
        $range.Value = $newStatus
    } else {
        # Fall back to an unconditional update if the text ever drifts.
        $range.Value = $newStatus
    }
}

# Re-fit the Status columns now that the text is shorter.
$overview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E
$overview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F
$zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C
$dede.Columns.Item(3).ColumnWidth = $newColumnWidth        # column C
